$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value (kept as Text to match the original inlineStr cell type)
$updates = @{
    "D2" = "286.62"
    "E2" = "4.39%"
    "D3" = "28.34"
    "E3" = "4.03%"
    "D4" = "4.920"
    "E4" = "0.86%"
    "D5" = "0.06549"
    "E5" = "2.47%"
    "E6" = "4.42%"
    "D7" = "1.364"
    "E7" = "14.75%"
    "D8" = "0.9141"
    "E8" = "4.23%"
    "E9" = "4.01%"
    "D10" = "0.06677"
    "E10" = "30.67%"
    "D11" = "0.07728"
    "E11" = "2.57%"
    "D12" = "0.02983"
    "E12" = "0.66%"
    "D13" = "0.08977"
    "D14" = "0.001594"
    "E14" = "2.09%"
    "D15" = "0.0006564"
    "E15" = "2.87%"
    "D16" = "0.006072"
    "E16" = "-1.83%"
    "D17" = "3.485"
    "E17" = "0.52%"
    "E18" = "2.53%"
    "D19" = "2.245"
    "E19" = "-1.69%"
    "E20" = "0.67%"
    "D22" = "3.978"
    "E22" = "1.40%"
    "D23" = "0.04460"
    "E23" = "1.11%"
    "D24" = "0.1520"
    "E24" = "10.13%"
    "D25" = "0.001186"
    "E25" = "0.66%"
    "D26" = "0.004337"
    "E26" = "12.56%"
    "E28" = "-1.80%"
    "D40" = "0.04155"
    "E40" = "0.70%"
    "D41" = "0.006904"
    "E41" = "1.56%"
    "D42" = "0.1412"
    "E42" = "20.31%"
    "E43" = "-6.08%"
    "D44" = "0.01244"
    "E44" = "8.35%"
    "D45" = "0.00005559"
    "E45" = "6.96%"
    "D46" = "1.561"
    "E46" = "-7.12%"
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
}

Write-Output "Updated $($updates.Count) cells"